# Aerodynamic features testing in progress
#
# The Xcg estimation method comparison blocks on the FUSELAGE and WING
# sheets had the TORENBEEK_1982 / SFORZA values accidentally swapped
# relative to their labels. Swap the numeric results back so each method
# label lines up with its own computed value.

$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet: rows 23 (TORENBEEK_1982) / 24 (SFORZA) ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
$fuseTorenbeek = $wsFuselage.Range("C23").Value2
$fuseSforza = $wsFuselage.Range("C24").Value2
$wsFuselage.Range("C23").Value = $fuseSforza
$wsFuselage.Range("C24").Value = $fuseTorenbeek

# --- WING sheet: two blocks, rows 23/24 and rows 27/28 ---
$wsWing = $wb.Worksheets.Item("WING")

$wingTorenbeek1 = $wsWing.Range("C23").Value2
$wingSforza1 = $wsWing.Range("C24").Value2
$wsWing.Range("C23").Value = $wingSforza1
$wsWing.Range("C24").Value = $wingTorenbeek1

$wingTorenbeek2 = $wsWing.Range("C27").Value2
$wingSforza2 = $wsWing.Range("C28").Value2
$wsWing.Range("C27").Value = $wingSforza2
$wsWing.Range("C28").Value = $wingTorenbeek2
